$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -81.9309
$ws.Range("B2").Value = -81.8674

$ws.Range("A3").Value = 25.6658
$ws.Range("B3").Value = 25.7247

$ws.Range("A4").Value = -81.338
$ws.Range("B4").Value = -81.4021

$ws.Range("A5").Value = 26.2119
$ws.Range("B5").Value = 26.1533
